$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("listName")

# Row 28: add a couple (B,C) that was previously missing
$ws.Range("B28").Value = "สมบูรณ์"
$ws.Range("C28").Value = "บริบูรณ์"

# Row 14: add a couple (B,C) that was previously missing
$ws.Range("B14").Value = "ปิยมาภรณ์"
$ws.Range("C14").Value = "รู้นันต๊ะ"

# Row 30: brand new row appended at the bottom
$ws.Range("A30").Value = "พี่มี่"
$ws.Range("B30").Value = "ทัศน์พล"
$ws.Range("C30").Value = "ผดุงโกเม็ด"

# Update the view: scroll so row 13 is at top, and select D25
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("D25").Select()
